$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; old column B (and its formatting)
# shifts to column C automatically.
$ws.Columns("B:B").Insert()

# Header for the new column.
$ws.Cells.Item(1, 2).Value = "deidentified"

# Rows that are NOT kept in the de-identified export (highlighted red
# in the original sheet) get a 0 flag; everything else gets 1.
$zeroRows = @(12, 13, 14)
for ($r = 2; $r -le 37; $r++) {
    if ($zeroRows -contains $r) {
        $ws.Cells.Item($r, 2).Value = 0
    } else {
        $ws.Cells.Item($r, 2).Value = 1
    }
}

# Clean, unformatted style for the whole new column (matches the plain
# "Normal" look used elsewhere in the workbook for blank/default cells).
$rng = $ws.Range("B1:B37")
$rng.ClearFormats()

# Column width for the new flag column (old column C keeps the width it
# already carried over from the pre-insert column B, so it's left alone).
$ws.Columns("B:B").ColumnWidth = 9.285714285714286

# The insert shifts the sheet's old "blank default" column definition
# (formerly column C) out to column D, where it's now unused -- drop its
# formatting so no stray <col> entry for it is written out.
$ws.Columns("D:D").ClearFormats()

# Restore the active selection.
$ws.Range("F6").Select()
